$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 45) with the latest bitcoin buy entry.
# Column A holds the date as literal text (e.g. "08/31/2025"), matching the
# other recently-appended rows in this sheet. Temporarily mark the cell as
# text so Excel does not auto-convert the string into a date serial number,
# then clear the formatting override so the cell keeps the sheet's default
# (unstyled) look, just like the other text-date cells above it.
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "08/31/2025"
$ws.Range("A45").ClearFormats()

$ws.Range("B45").Value = 0.0004578099999999995
$ws.Range("C45").Value = 109215.613464101
$ws.Range("D45").Value = 50
